# Revert capacity chart to show kilowatts (instead of watts) on the y-axis.
#
# The underlying "Solar" column (column E, rows 20-26 / years 2018-2024)
# was previously stored in watts; convert the stored values (and the
# number format applied to the whole data body) to kilowatts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Convert the raw data values in column E (Solar) from watts to
#    kilowatts for the affected rows (2018-2024).
$ws.Range("E20").Value = 6
$ws.Range("E21").Value = 5
$ws.Range("E22").Value = 3.1
$ws.Range("E23").Value = 30.02
$ws.Range("E24").Value = 9.24
$ws.Range("E25").Value = 56.03
$ws.Range("E26").Value = 44.22

# 2) The shared number format used across the whole data table (B2:G26)
#    gains one decimal place so fractional kilowatt values display
#    correctly.
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# 3) Update the chart's value-axis title and number format so the axis
#    reads "Kilowatts (kW)" with a plain thousands format instead of the
#    "Watts" title with a ">=1000 => K" abbreviation format.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$valueAxis = $chart.Axes(2)
$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"
